# Weekly update: insert two new report rows at the top of the data block
# (row 837), pushing the existing rows down by two. This mirrors the
# "Fruta / hortaliza, semanal" commit: a new week's Pimiento (Zafiro rojo /
# Zafiro verde, Región de Arica y Parinacota) entries are prepended.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 837 (shifts rows 837:909 down to 839:911).
$ws.Rows.Item(837).Insert()
$ws.Rows.Item(837).Insert()

# New row 837: Zafiro rojo
$ws.Range("A837").Value = 5
$ws.Range("B837").Value = "Macroferia Regional de Talca"
$ws.Range("C837").Value = "Maule"
$ws.Range("D837").Value = 45106
$ws.Range("E837").Value = 7
$ws.Range("F837").Value = 100112002
$ws.Range("G837").Value = "Pimiento"
$ws.Range("H837").Value = "Zafiro rojo"
$ws.Range("I837").Value = "Primera"
$ws.Range("J837").Value = 300
$ws.Range("K837").Value = 10000
$ws.Range("L837").Value = 10000
$ws.Range("M837").Value = 10000
$ws.Range("N837").Value = "$/caja 15 kilos"
$ws.Range("O837").Value = "Región de Arica y Parinacota"
$ws.Range("P837").Value = 667
$ws.Range("Q837").Value = 15
$ws.Range("R837").Value = "Hortaliza"

# New row 838: Zafiro verde
$ws.Range("A838").Value = 5
$ws.Range("B838").Value = "Macroferia Regional de Talca"
$ws.Range("C838").Value = "Maule"
$ws.Range("D838").Value = 45106
$ws.Range("E838").Value = 7
$ws.Range("F838").Value = 100112002
$ws.Range("G838").Value = "Pimiento"
$ws.Range("H838").Value = "Zafiro verde"
$ws.Range("I838").Value = "Primera"
$ws.Range("J838").Value = 300
$ws.Range("K838").Value = 10000
$ws.Range("L838").Value = 10000
$ws.Range("M838").Value = 10000
$ws.Range("N838").Value = "$/caja 15 kilos"
$ws.Range("O838").Value = "Región de Arica y Parinacota"
$ws.Range("P838").Value = 667
$ws.Range("Q838").Value = 15
$ws.Range("R838").Value = "Hortaliza"
